# Update metrics_table.xlsx with recomputed SHAP/RF optimisation scores
# for the "Clinical ABC / top / SHAP" metrics table (rows 4-8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6039309056956117
$ws.Range("C4").Value = 0.05182547567860592
$ws.Range("I4").Value = 0.6818063895514784
$ws.Range("J4").Value = 0.7029454022988506
$ws.Range("O4").Value = 0.6849653565129756
$ws.Range("P4").Value = 0.0320721661651833
$ws.Range("B5").Value = 0.6151695512044775
$ws.Range("C5").Value = 0.008033281799195885
$ws.Range("D5").Value = 0.6040268456375839
$ws.Range("E5").Value = 0.7679180887372016
$ws.Range("F5").Value = 0.4433962264150944
$ws.Range("G5").Value = 0.4455445544554456
$ws.Range("H5").Value = 0.9375
$ws.Range("I5").Value = 0.480220032029797
$ws.Range("J5").Value = 0.4752155172413793
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 56
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 45
$ws.Range("O5").Value = 0.5984282354917275
$ws.Range("P5").Value = 0.03076535305184156
$ws.Range("B6").Value = 0.5908779683957377
$ws.Range("C6").Value = 0.1012384708695192
$ws.Range("D6").Value = 0.6326530612244898
$ws.Range("E6").Value = 0.6404958677685951
$ws.Range("F6").Value = 0.660377358490566
$ws.Range("G6").Value = 0.62
$ws.Range("H6").Value = 0.6458333333333334
$ws.Range("I6").Value = 0.6060940472510277
$ws.Range("J6").Value = 0.6663074712643678
$ws.Range("M6").Value = 17
$ws.Range("N6").Value = 31
$ws.Range("O6").Value = 0.6542584194171497
$ws.Range("P6").Value = 0.0559439279691001
$ws.Range("B7").Value = 0.5118928236061546
$ws.Range("C7").Value = 0.05444872521682439
$ws.Range("D7").Value = 0.4444444444444445
$ws.Range("E7").Value = 0.4273504273504274
$ws.Range("F7").Value = 0.5283018867924528
$ws.Range("G7").Value = 0.4761904761904762
$ws.Range("I7").Value = 0.4625636418089248
$ws.Range("J7").Value = 0.5186781609195403
$ws.Range("K7").Value = 36
$ws.Range("L7").Value = 22
$ws.Range("O7").Value = 0.5653174603174603
$ws.Range("P7").Value = 0.03634023286724132
$ws.Range("B8").Value = 0.5802362244131926
$ws.Range("C8").Value = 0.06104593199208728
$ws.Range("D8").Value = 0.4819277108433734
$ws.Range("E8").Value = 0.4405286343612335
$ws.Range("F8").Value = 0.5943396226415094
$ws.Range("G8").Value = 0.5714285714285714
$ws.Range("H8").Value = 0.4166666666666667
$ws.Range("I8").Value = 0.5796564106033033
$ws.Range("J8").Value = 0.6131465517241379
$ws.Range("K8").Value = 43
$ws.Range("L8").Value = 15
$ws.Range("M8").Value = 28
$ws.Range("N8").Value = 22
$ws.Range("O8").Value = 0.6701383639875703
$ws.Range("P8").Value = 0.04446710256275315
